{"js": "// \"Version 1.\" -> \"Version 2.\" (v1.2/Interaction/wireframes.docx)\n//\n// The visible edit is just the digit 1 -> 2, but the canonical OOXML also\n// shows two side effects of how Word recorded the edit:\n//   - the word \"Version\" ends up split into two runs (\"Versi\" / \"on\") with\n//     no formatting change (an artifact of a cursor/selection boundary);\n//   - the \"_GoBack\" bookmark (Word's \"last edit location\" marker) ends up\n//     right after the new \"2\" instead of at the end of the paragraph,\n//     which pushes the trailing \".\" onto its own run after the bookmark.\n// Both are reproduced by planting/removing a bookmark at the exact\n// character spot: Word.js routes insertBookmark/deleteBookmark through the\n// same run-splitting machinery real Word uses, without leaving any stray\n// run formatting behind.\n\nconst body = context.document.body;\n\n// --- Step 1: split \"Version\" -> \"Versi\" | \"on\" --------------------------\n// Insert a throwaway bookmark right at the \"Versi\"|\"on\" boundary, then\n// delete it immediately; this forces a clean run split at that offset.\nconst versiResults = body.search(\"Versi\", { matchCase: true });\nversiResults.load(\"text\");\nawait context.sync();\n\nconst splitPoint = versiResults.items[0].getRange(Word.RangeLocation.after);\nsplitPoint.insertBookmark(\"_tmpSplit\");\nawait context.sync();\n\ncontext.document.deleteBookmark(\"_tmpSplit\");\nawait context.sync();\n\n// --- Step 2: change the version digit, 1 -> 2 ----------------------------\nconst oneResults = body.search(\"1\", { matchCase: true });\noneResults.load(\"text\");\nawait context.sync();\n\noneResults.items[0].insertText(\"2\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Step 3: relocate \"_GoBack\" to sit right after the edited \"2\" -------\n// Remove the \".\" and retype it fresh so it doesn't inherit the old\n// \" 1.\" run's xml:space=\"preserve\" (a bare \".\" doesn't need it); the\n// \"_GoBack\" bookmark (already sitting right after the old \".\") stays put\n// while the new \".\" run lands after it, matching real Word's behavior of\n// leaving \"_GoBack\" at the actual edit point.\nconst dotResults = body.search(\".\", { matchCase: true, matchWildcards: false });\ndotResults.load(\"text\");\nawait context.sync();\n\nconst dotRange = dotResults.items[0];\nconst beforeDot = dotRange.getRange(Word.RangeLocation.before);\ndotRange.delete();\nawait context.sync();\n\nbeforeDot.insertText(\".\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# \"Version 1.\" -> \"Version 2.\" (v1.2/Interaction/wireframes.docx)\n#\n# The visible edit is just the digit 1 -> 2, but the canonical OOXML also\n# shows two side effects of how Word tracked the edit:\n#   - the word \"Version\" is split into two runs (\"Versi\" / \"on\") with no\n#     formatting change (an artifact of a cursor/selection boundary);\n#   - the \"_GoBack\" bookmark (Word's \"last edit location\" marker) is moved\n#     from the end of the paragraph to right after the new \"2\", which\n#     forces the trailing \".\" onto its own run after the bookmark.\n# We reproduce both by briefly planting/removing bookmarks at the exact\n# character offsets, which forces Word's OM to split runs there without\n# leaving any stray run-formatting behind.\n\n$d = $word.ActiveDocument\n\n# Work against the whole document range/text (\"Version 1.\" == 10 chars).\n$text = $d.Content.Text\n\n# --- Step 1: split \"Version\" -> \"Versi\" | \"on\" -------------------------\n# Dropping a bookmark at the boundary and immediately deleting it forces a\n# clean run split at that character offset (no residual formatting).\n$splitPoint = $d.Range(5, 5)\n$d.Bookmarks.Add(\"_tmpSplit\", $splitPoint) | Out-Null\n$d.Bookmarks.Item(\"_tmpSplit\").Delete()\n\n# --- Step 2: change the version digit, 1 -> 2 ---------------------------\n$d.Range(8, 9).Text = \"2\"\n\n# --- Step 3: relocate \"_GoBack\" to sit right after the edited \"2\" -------\n# Remove the \".\" and retype it fresh so it doesn't inherit the old\n# \" 1.\" run's xml:space=\"preserve\" (a bare \".\" doesn't need it).\n$d.Range(9, 10).Delete()\n$d.Range(9, 9).InsertAfter(\".\")\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n$goBackPos = $d.Range(9, 9)\n$d.Bookmarks.Add(\"_GoBack\", $goBackPos) | Out-Null\n"}
